$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "App()" row (old row 2); remaining rows shift up.
$ws.Range("A2:K2").EntireRow.Delete()

# Row 2 (was row 3): main(String[] args)
$ws.Cells.Item(2, 1).Value = 1.0
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 4.0
$ws.Cells.Item(2, 9).Value = 4.0
$ws.Cells.Item(2, 10).Value = 1.0

# Row 3 (was row 4): memoryRecall()
$ws.Cells.Item(3, 1).Value = 2.0
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 4.0
$ws.Cells.Item(3, 9).Value = 5.0
$ws.Cells.Item(3, 10).Value = 2.0

# Row 4 (was row 5): memoryReca()
$ws.Cells.Item(4, 1).Value = 3.0
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 4.0
$ws.Cells.Item(4, 9).Value = 5.0
$ws.Cells.Item(4, 10).Value = 1.0
